# Implement date and "by" field auto-population UI changes:
# Insert two new data rows (for two new shipments) above the previous
# last three rows, pushing the existing data down.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at row 6, shifting existing rows 6-8 down to 8-10.
$ws.Range("A6:A7").EntireRow.Insert()

# Helper to set a text value in a cell while forcing string storage
# (even for numeric-looking values) and keeping the default "Normal" style.
function Set-TextCell($row, $col, $text) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

# New row 6: Mars / JS01 / AC / 23 / Space Weed / SW420 / kLi3pJdPSNnEaSk0
Set-TextCell 6 1 "Mars"
Set-TextCell 6 2 "JS01"
Set-TextCell 6 3 "AC"
Set-TextCell 6 4 "23"
Set-TextCell 6 5 "Space Weed"
Set-TextCell 6 6 "SW420"
Set-TextCell 6 7 "kLi3pJdPSNnEaSk0"

# New row 7: TO4 / MAIN / AC / 2 / test / test / l8fWIfXMgODcFPYQ
Set-TextCell 7 1 "TO4"
Set-TextCell 7 2 "MAIN"
Set-TextCell 7 3 "AC"
Set-TextCell 7 4 "2"
Set-TextCell 7 5 "test"
Set-TextCell 7 6 "test"
Set-TextCell 7 7 "l8fWIfXMgODcFPYQ"
